$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 held only B13/C13 ("519033 - Carlos Yujiro Shigue", no A label).
# Deleting it shifts rows 14-25 up to 13-24; row heights and the sheet dimension follow
# automatically with the shift.
$ws.Rows(13).Delete()

# After the shift, the B/C content of several surviving rows must be updated to the
# new canonical text for that row.
$ws.Range("B10").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C10").Value = '519033 - Carlos Yujiro Shigue'

$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

$ws.Range("B18").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C18").Value = '519033 - Carlos Yujiro Shigue'

$ws.Range("B19").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Range("C19").Value = 'Aulas expositivas, seminários e exercícios comentados.'

$ws.Range("B20").Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$ws.Range("C20").Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'

$ws.Range("B21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'

# Row 15 needs the literal text "01/01/2012" (same text already used verbatim in row 8).
# Assigning it with .Value would get auto-parsed into a date serial number, so copy the
# existing text cell instead, which preserves it as a plain shared string.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

